$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (Förändrad) for rows 2-5: bump date serial from 45174 to 45175 (2023-09-05 -> 2023-09-06)
$ws.Range("C2").Value = 45175
$ws.Range("C3").Value = 45175
$ws.Range("C4").Value = 45175
$ws.Range("C5").Value = 45175
